# Update leaderboard worksheets to reflect the 27.03.2025 20:00 standings.

$wb = $excel.ActiveWorkbook

# --- Sheet "leaderboard2" (Cobblemons caught) ---
$ws1 = $wb.Worksheets.Item("leaderboard2")
$ws1.Range("D3").Value = 261
$ws1.Range("D4").Value = 84
$ws1.Range("B13").Value = "Dernière update le 27.03.25 à 20:00"

# --- Sheet "leaderboard3" (Shiny Cobblemons caught) - 1st/2nd place swapped ---
$ws2 = $wb.Worksheets.Item("leaderboard3")
$ws2.Range("C3").Value = "ArtyumsM"
$ws2.Range("D3").Value = 19
$ws2.Range("C4").Value = "BKZRackham"
$ws2.Range("D4").Value = 19
$ws2.Range("B13").Value = "Dernière update le 27.03.25 à 20:00"
